$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column AF header
$ws.Range("AF1").Value = "21-jul"

# New column AF values, row 2-11 (match formatting of column AE)
$values = @(11, 15, 9, 14, 15, 13, 18, 15, 23, 22)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 31).Copy()
    $ws.Cells.Item($row, 32).PasteSpecial(-4122)
    $ws.Cells.Item($row, 32).Value = $values[$i]
}

# Update selection to match the new active cell
$ws.Range("AF12").Select()
